$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - dnn_n51_transpiled.qasm
$ws.Range("B5").Value = 0.003659274240275831
$ws.Range("C5").Value = 0.002797244170942612
$ws.Range("D5").Value = 269
$ws.Range("E5").Value = 252

# Row 7 - sqrt18.qasm
$ws.Range("B7").Value = 0.00003829878473929803
$ws.Range("C7").Value = 0.00003030098235237177
$ws.Range("D7").Value = 898
$ws.Range("E7").Value = 856

# Row 8 - dnn_n33_transpiled.qasm
$ws.Range("B8").Value = 0.0707183695196887
$ws.Range("C8").Value = 0.06255269054031519
$ws.Range("D8").Value = 170
$ws.Range("E8").Value = 158

# Row 9 - qft_n18.qasm
$ws.Range("B9").Value = 0.03262570547925792
$ws.Range("C9").Value = 0.02984051258584217
$ws.Range("D9").Value = 300
$ws.Range("E9").Value = 284

# Row 10 - DNN16.qasm
$ws.Range("C10").Value = 0.3102413784182018
$ws.Range("D10").Value = 87

# Row 11 - QV_32.qasm
$ws.Range("B11").Value = 0.0000000001302193996108339
$ws.Range("C11").Value = 0.00000000007738876542690178
$ws.Range("D11").Value = 1447
$ws.Range("E11").Value = 1403
$ws.Range("G11").Value = 1473

# Row 13 - hhl_n7.qasm
$ws.Range("B13").Value = 0.507471706605126
$ws.Range("E13").Value = 81

# Row 14 - qaoa_n6_transpiled.qasm
$ws.Range("C14").Value = 0.7589543307708074
$ws.Range("D14").Value = 32

# Row 16 - dder_n10_transpiled.qasm
$ws.Range("C16").Value = 0.6250617265957956
$ws.Range("D16").Value = 64
